$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
try {
    $d.Delete()
    Write-Output "DELETE OK"
} catch {
    Write-Output ("ERR: " + $_.Exception.Message)
}
Write-Output ("Designs.Count=" + $p.Designs.Count)
